# Updated TestData for Portugal Market
# - Add a new "Portugal" worksheet (copied from the "Czech" template sheet,
#   which shares its row/column layout) as the last tab, with market name
#   "Portugal Market" and JIRA ref "NGC-3479/T1756".
# - Make "Portugal" the active tab (previously "Swiss" / index 3 was active).
# - Fix up the "Germany" sheet selection so it reflects the used range
#   instead of a full-sheet selection.

$wb = $excel.ActiveWorkbook

# --- Fix Germany sheet selection (was selecting the whole sheet) ---
$germany = $wb.Worksheets.Item("Germany")
$germany.Range("A1:D21").Select() | Out-Null

# --- Create the new Portugal sheet from the Czech sheet template ---
$czech = $wb.Worksheets.Item("Czech")
$swiss = $wb.Worksheets.Item("Swiss")
$czech.Copy($null, $swiss)

$portugal = $wb.Worksheets.Item($wb.Worksheets.Count)
$portugal.Name = "Portugal"

# --- Market-specific cell content ---
$portugal.Range("B2").Value = "Portugal Market"
$portugal.Range("B4").Value = "NGC-3479/T1756"

# --- Row heights for the (now taller/wrapped) header rows ---
$portugal.Rows.Item(3).RowHeight = 28.8
$portugal.Rows.Item(4).RowHeight = 28.8
$portugal.Rows.Item(5).RowHeight = 28.8

# --- Column widths for the new sheet's layout ---
$portugal.Columns.Item(1).ColumnWidth = 27.72
$portugal.Columns.Item(2).ColumnWidth = 18.05
$portugal.Columns.Item(3).ColumnWidth = 11.5
$portugal.Columns.Item(4).ColumnWidth = 14.17

# --- Selection + activate (makes it the selected/active tab) ---
$portugal.Range("O17").Select() | Out-Null
$portugal.Activate() | Out-Null
